$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng24 = $ws.Range("J24")
$rng16 = $ws.Range("J16")
Write-Output ("row24 bottom ColorIndex=" + $rng24.Borders(9).ColorIndex() + " ThemeColor=" + $rng24.Borders(9).ThemeColor())
Write-Output ("row16 bottom ColorIndex=" + $rng16.Borders(9).ColorIndex() + " ThemeColor=" + $rng16.Borders(9).ThemeColor())
